$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update F column values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 80
$ws1.Range("F7").Value = 97
$ws1.Range("F13").Value = 2417
$ws1.Range("F14").Value = 32
$ws1.Range("F18").Value = 529
$ws1.Range("F19").Value = 564
$ws1.Range("F24").Value = 1972
$ws1.Range("F25").Value = 4100
$ws1.Range("F28").Value = 1197
$ws1.Range("F29").Value = 233
$ws1.Range("F30").Value = 2104
$ws1.Range("F35").Value = 290
$ws1.Range("F38").Value = 709

# Sheet "全部类型" (fourth sheet) - update F column values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 80
$ws4.Range("F7").Value = 97
$ws4.Range("F13").Value = 2417
$ws4.Range("F14").Value = 32
$ws4.Range("F19").Value = 529
$ws4.Range("F20").Value = 564
$ws4.Range("F25").Value = 1972
$ws4.Range("F26").Value = 4100
$ws4.Range("F29").Value = 1197
$ws4.Range("F30").Value = 233
$ws4.Range("F31").Value = 2104
$ws4.Range("F36").Value = 290
$ws4.Range("F39").Value = 709

$wb.Save()
